$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.960.84"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").Value = "2.790.02"
$ws.Range("E3").Value = "  -1.80%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "358.32"
$ws.Range("E5").Value = "  -0.98%  "
$ws.Range("D6").Value = "109.67"
$ws.Range("E6").Value = "  -3.48%  "
$ws.Range("D7").Value = "0.557"
$ws.Range("E7").Value = "  -1.82%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "0.589"
$ws.Range("E9").Value = "  -2.25%  "
$ws.Range("D10").Value = "39.78"
$ws.Range("E10").Value = "  -4.56%  "
$ws.Range("D11").Value = "0.0846"
$ws.Range("E11").Value = "  -1.66%  "
$ws.Range("D13").Value = "19.41"
$ws.Range("E13").Value = "  -3.20%  "
$ws.Range("D14").Value = "7.54"
$ws.Range("E14").Value = "  -3.08%  "
$ws.Range("D15").Value = "3.227.15"
$ws.Range("E15").Value = "  -1.91%  "
$ws.Range("D16").Value = "2.777.48"
$ws.Range("E16").Value = "  -2.18%  "
$ws.Range("D17").Value = "0.936"
$ws.Range("E17").Value = "  +3.32%  "
$ws.Range("D18").Value = "51.922.91"
$ws.Range("E18").Value = "  -0.02%  "
$ws.Range("D19").Value = "7.45"
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("E20").Value = "  -2.33%  "
$ws.Range("D21").Value = "13.07"
$ws.Range("E21").Value = "  -3.52%  "
$ws.Range("D22").Value = "0.0₃0975"
$ws.Range("E22").Value = "  -1.76%  "
$ws.Range("D23").Value = "70.17"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").Value = "269.68"
$ws.Range("E24").Value = "  +0.96%  "
$ws.Range("D25").Value = "2.74"
$ws.Range("E25").Value = "  -3.47%  "
$ws.Range("D26").Value = "26.51"
$ws.Range("E26").Value = "  -2.23%  "
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("D28").Value = "0.163"
$ws.Range("E28").Value = "  +16.16%  "
$ws.Range("D29").Value = "10.23"
$ws.Range("E29").Value = "  -1.70%  "
$ws.Range("E30").Value = "  -1.64%  "
$ws.Range("D31").Value = "0.0470"
$ws.Range("E31").Value = "  +5.30%  "
$ws.Range("D32").Value = "51.89"
$ws.Range("E32").Value = "  -3.63%  "
$ws.Range("D33").Value = "33.66"
$ws.Range("E33").Value = "  -1.24%  "
$ws.Range("D34").Value = "5.72"
$ws.Range("E34").Value = "  -2.89%  "
$ws.Range("E35").Value = "  -0.29%  "
$ws.Range("D36").Value = "5.19"
$ws.Range("E36").Value = "  -2.35%  "
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("D38").Value = "18.78"
$ws.Range("E38").Value = "  +2.46%  "
$ws.Range("E39").Value = "  -2.97%  "
$ws.Range("D40").Value = "1.99"
$ws.Range("E40").Value = "  -4.24%  "
$ws.Range("D41").Value = "2.55"
$ws.Range("E41").Value = "  -0.51%  "
$ws.Range("E42").Value = "  -1.70%  "
$ws.Range("E43").Value = "  -0.54%  "
$ws.Range("D44").Value = "119.39"
$ws.Range("E44").Value = "  -6.82%  "
$ws.Range("D45").Value = "21.68"
$ws.Range("E45").Value = "  -10.37%  "
$ws.Range("D46").Value = "2.081.23"
$ws.Range("E46").Value = "  -1.89%  "
$ws.Range("D47").Value = "3.23"
$ws.Range("E47").Value = "  -4.33%  "
$ws.Range("D48").Value = "2.21"
$ws.Range("E48").Value = "  -1.80%  "
$ws.Range("D49").Value = "5.81"
$ws.Range("E49").Value = "  -0.32%  "
$ws.Range("D50").Value = "0.953"
$ws.Range("E50").Value = "  -4.55%  "
$ws.Range("D51").Value = "8.87"
$ws.Range("E51").Value = "  -1.79%  "
